$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to Text format so numeric-looking price strings
# (e.g. "1.00", "606.60") keep their exact text representation
# instead of being auto-converted to numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '66.036.64'
$ws.Range('E2').Value = '  +0.09%  '
$ws.Range('D3').Value = '3.174.25'
$ws.Range('E3').Value = '  -0.99%  '
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').Value = '606.60'
$ws.Range('E5').Value = '  +0.82%  '
$ws.Range('D6').Value = '154.63'
$ws.Range('E6').Value = '  +0.74%  '
$ws.Range('D7').Value = '1.00'
$ws.Range('E7').Value = '  +0.02%  '
$ws.Range('D8').Value = '3.171.97'
$ws.Range('E8').Value = '  -0.94%  '
$ws.Range('D9').Value = '0.545'
$ws.Range('E9').Value = '  +2.28%  '
$ws.Range('E10').Value = '  -1.07%  '
$ws.Range('E11').Value = '  -6.66%  '
$ws.Range('E12').Value = '  +1.13%  '
$ws.Range('E13').Value = '  -1.74%  '
$ws.Range('D14').Value = '38.23'
$ws.Range('E14').Value = '  -3.18%  '
$ws.Range('D15').Value = '3.695.17'
$ws.Range('E15').Value = '  -0.96%  '
$ws.Range('D16').Value = '66.122.54'
$ws.Range('E16').Value = '  +0.15%  '
$ws.Range('D18').Value = '3.177.29'
$ws.Range('E18').Value = '  -0.49%  '
$ws.Range('E19').Value = '  +1.08%  '
$ws.Range('D20').Value = '508.56'
$ws.Range('E20').Value = '  -0.63%  '
$ws.Range('D21').Value = '15.34'
$ws.Range('E21').Value = '  -0.71%  '
$ws.Range('D22').Value = '0.728'
$ws.Range('E22').Value = '  -1.61%  '
$ws.Range('D23').Value = '8.00'
$ws.Range('E23').Value = '  -2.21%  '
$ws.Range('D24').Value = '14.76'
$ws.Range('E24').Value = '  -4.14%  '
$ws.Range('D25').Value = '84.37'
$ws.Range('E25').Value = '  -0.65%  '
$ws.Range('E26').Value = '  +0.10%  '
$ws.Range('D28').Value = '9.15'
$ws.Range('E28').Value = '  -2.11%  '
$ws.Range('E29').Value = '  +4.27%  '
$ws.Range('E30').Value = '  +4.69%  '
$ws.Range('E31').Value = '  +4.20%  '
$ws.Range('D32').Value = '27.92'
$ws.Range('E32').Value = '  -0.67%  '
$ws.Range('E33').Value = '  +0.27%  '
$ws.Range('D34').Value = '1.18'
$ws.Range('E34').Value = '  -3.25%  '
$ws.Range('D35').Value = '6.49'
$ws.Range('E35').Value = '  -1.19%  '
$ws.Range('D36').Value = '506.82'
$ws.Range('E36').Value = '  +3.91%  '
$ws.Range('D37').Value = '55.33'
$ws.Range('E37').Value = '  +0.62%  '
$ws.Range('D38').Value = '0.0877'
$ws.Range('E38').Value = '  -3.14%  '
$ws.Range('D39').Value = '0.0418'
$ws.Range('E39').Value = '  -0.69%  '
$ws.Range('E40').Value = '  +6.04%  '
$ws.Range('D41').Value = '8.76'
$ws.Range('E41').Value = '  -1.70%  '
$ws.Range('D42').Value = '0.0₃0684'
$ws.Range('E42').Value = '  +5.76%  '
$ws.Range('D43').Value = '2.84'
$ws.Range('E43').Value = '  -3.73%  '
$ws.Range('D44').Value = '0.296'
$ws.Range('E44').Value = '  -2.72%  '
$ws.Range('E45').Value = '  -0.36%  '
$ws.Range('D46').Value = '2.831.26'
$ws.Range('E46').Value = '  -4.16%  '
$ws.Range('D47').Value = '28.02'
$ws.Range('E47').Value = '  -2.15%  '
$ws.Range('D49').Value = '2.35'
$ws.Range('E49').Value = '  +1.84%  '
$ws.Range('E50').Value = '  +0.16%  '
$ws.Range('E51').Value = '  +3.00%  '

# Restore the default cell style on column D so no stray number-format
# style is left attached to the cells (matches original formatting).
$ws.Range("D2:D51").Style = "Normal"
